# ---------------------------------------------------------------------------
# New mangrove survey and WQ data
# Appends rows 28-37 to the "raw_data" sheet (dates 2022-04-07 .. 2022-04-16),
# refreshes a couple of workbook/sheet view settings, and re-creates the two
# new shared strings ("C3" and "<0.01") plus the two new cell styles that the
# authoring session introduced (a plain date style for column B, and a
# black-font style used on some of the J/K "Nitrate_raw"/"Nitrate" cells).
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("raw_data")

# ---------------------------------------------------------------------------
# 0. Row data. Columns: A Site, B Date(serial), C Type, D ph, E Conductivity,
#    F (unused/blank), G TDS, H Turbidity, I DO, J Nitrate_raw, K Nitrate,
#    L Phosphate_raw, M Phosphate, N Temperature.
#    J/K carry either the literal string "<0.01" or a plain number.
# ---------------------------------------------------------------------------
$rows = @(
    @(28, "T1", 44658, "ND", 8.1300000000000008, 34400, 30450, 0.88,               6.51, "<0.01", 0,    0.33,               0.33,               26.4, $false),
    @(29, "T2", 44659, "ND", 8.1999999999999993, 34200, 30850, 1.0900000000000001, 6.52, "<0.01", 0,    0.26,               0.26,               26.5, $false),
    @(30, "N1", 44660, "ND", 8.1999999999999993, 33900, 30600, 0.89,               6.63, 0.01,     0.01, 0.28999999999999998, 0.28999999999999998, 26.4, $false),
    @(31, "N2", 44661, "ND", 8.18,               33800, 30880, 0.63,               6.46, "<0.01", 0,    0.25,               0.25,               26.5, $true),
    @(32, "T4", 44662, "ND", 8.17,               33700, 30320, 0.7,                6.54, "<0.01", 0,    0.23,               0.23,               26.4, $true),
    @(33, "T5", 44663, "ND", 8.2200000000000006, 34100, 30550, 0.57999999999999996,7.01, "<0.01", 0,    0.2,                0.2,                26.5, $true),
    @(34, "T6", 44664, "ND", 8.19,               34000, 30490, 0.52,               6.84, 0.02,     0.02, 0.21,               0.21,               26.9, $false),
    @(35, "C1", 44665, "ND", 8.19,               34300, 30320, 0.62,               6.56, "<0.01", 0,    0.2,                0.2,                26.8, $true),
    @(36, "C3", 44666, "ND", 8.18,               33700, 30550, 0.65,               6.6,  "<0.01", 0,    0.25,               0.25,               26.4, $true),
    @(37, "T3", 44667, "ND", 8.17,               34400, 30810, 0.72,               6.56, "<0.01", 0,    0.26,               0.26,               26.6, $true)
)

# ---------------------------------------------------------------------------
# 1. Seed the shared-string table in the same order the original file has it:
#    "C3" (site id, first used on row 36) must land before "<0.01" (first
#    used on row 28) so the new entries come out as index 59 and 60.
# ---------------------------------------------------------------------------
$ws.Range("A36").Value = "C3"
$ws.Range("J28").Value = "<0.01"

# ---------------------------------------------------------------------------
# 2. Write every cell value for the new rows.
# ---------------------------------------------------------------------------
foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A$rowNum").Value = $r[1]
    $ws.Range("B$rowNum").Value = $r[2]
    $ws.Range("C$rowNum").Value = $r[3]
    $ws.Range("D$rowNum").Value = $r[4]
    $ws.Range("E$rowNum").Value = $r[5]
    $ws.Range("G$rowNum").Value = $r[6]
    $ws.Range("H$rowNum").Value = $r[7]
    $ws.Range("I$rowNum").Value = $r[8]
    $ws.Range("J$rowNum").Value = $r[9]
    $ws.Range("K$rowNum").Value = $r[10]
    $ws.Range("L$rowNum").Value = $r[11]
    $ws.Range("M$rowNum").Value = $r[12]
    $ws.Range("N$rowNum").Value = $r[13]
}

# ---------------------------------------------------------------------------
# 3. Formatting.
#    3a. Base data-row style (font 10pt, right/vcenter align - this is the
#        same style already used by rows 12-27) for every populated column
#        except B (date).
# ---------------------------------------------------------------------------
foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("A20").Copy()
    $ws.Range("A" + $rowNum + ":E" + $rowNum).PasteSpecial(-4122)
    $ws.Range("A20").Copy()
    $ws.Range("G" + $rowNum + ":N" + $rowNum).PasteSpecial(-4122)
}

# 3b. Column B: plain date style (no custom font/alignment), m/d/yyyy.
foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Range("ZZ1").Copy()
    $ws.Range("B$rowNum").PasteSpecial(-4122)
    $ws.Range("B$rowNum").NumberFormat = "m/d/yyyy"
}

# 3c. J/K black-font style on the rows that got it during the authoring
#     session (N2, T4, T5, C1, C3, T3 rows).
foreach ($r in $rows) {
    if ($r[14]) {
        $rowNum = $r[0]
        $ws.Range("A20").Copy()
        $ws.Range("J" + $rowNum + ":K" + $rowNum).PasteSpecial(-4122)
        $ws.Range("J" + $rowNum + ":K" + $rowNum).Font.Color = 0
    }
}

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4. Sheet/window view cosmetics: zoom + split panes around D18, matching the
#    reviewer's on-screen state when the new rows were entered.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.Zoom = 167
$excel.ActiveWindow.SplitColumn = 2
$excel.ActiveWindow.SplitRow = 4
$ws.Range("E34").Select()

Write-Output "done"
